$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "sbPj0AEACAAJ"
$ws.Range("B9").Value = "Juego de tronos"
$ws.Range("C9").Value = "Desconocido"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2018"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "George R.R Martin"
